# fix handling empty rows
# Inserts blank separator rows into the "transitionstates", "scenarios" and
# "businessevents" tables/sheets, resizes the backing Excel Tables to match
# the new extents, and restores the view state (selection / active sheet /
# frozen-pane scroll position) to what the user ended up looking at.

$wb = $excel.ActiveWorkbook

$wsTransitionStates = $wb.Worksheets.Item(1)   # "transitionstates"
$wsScenarios        = $wb.Worksheets.Item(2)   # "scenarios"
$wsBusinessEvents   = $wb.Worksheets.Item(3)   # "businessevents"

# ---------------------------------------------------------------------------
# 1) transitionstates: a blank row before every data row (rows 3 and 4 each
#    get pushed down, with one gap row inserted ahead of each).
# ---------------------------------------------------------------------------
$wsTransitionStates.Rows.Item(3).Insert()
$wsTransitionStates.Rows.Item(5).Insert()

$wsTransitionStates.ListObjects.Item(1).Resize($wsTransitionStates.Range("A1:C6"))

# ---------------------------------------------------------------------------
# 2) scenarios: same pattern - a blank row before every data row.
# ---------------------------------------------------------------------------
$wsScenarios.Rows.Item(3).Insert()
$wsScenarios.Rows.Item(5).Insert()
$wsScenarios.Rows.Item(7).Insert()
$wsScenarios.Rows.Item(9).Insert()

$wsScenarios.ListObjects.Item(1).Resize($wsScenarios.Range("A1:C10"))

# ---------------------------------------------------------------------------
# 3) businessevents: a blank row only between groups of rows that share the
#    same ScenarioKey (column C) - i.e. before rows that started a new group.
# ---------------------------------------------------------------------------
$wsBusinessEvents.Rows.Item(5).Insert()
$wsBusinessEvents.Rows.Item(9).Insert()
$wsBusinessEvents.Rows.Item(12).Insert()
$wsBusinessEvents.Rows.Item(15).Insert()

# Give the newly inserted blank rows the same "separator row" look used
# elsewhere in this sheet (col A keeps the normal key style, col B/C borrow
# the alternate styles already present on the sheet) rather than whatever
# default formatting Insert() picked up from the rows above them.
foreach ($r in 5, 9, 12, 15) {
    $wsBusinessEvents.Range("A2").Copy()
    $wsBusinessEvents.Range("A$r").PasteSpecial(-4122)
    $wsBusinessEvents.Range("B11").Copy()
    $wsBusinessEvents.Range("B$r").PasteSpecial(-4122)
    $wsBusinessEvents.Range("B2").Copy()
    $wsBusinessEvents.Range("C$r").PasteSpecial(-4122)
}

$wsBusinessEvents.ListObjects.Item(1).Resize($wsBusinessEvents.Range("A1:C17"))

# ---------------------------------------------------------------------------
# 4) View state: the user re-selected A3 on transitionstates, scrolled /
#    selected A5 on businessevents (frozen header stays at row 1, but the
#    window is now scrolled so row 2 is the first visible row under it), and
#    finally ended up on scenarios with A3 selected (making it the active
#    tab).
# ---------------------------------------------------------------------------
$wsTransitionStates.Range("A3").Select()

$wsBusinessEvents.Select()
$excel.ActiveWindow.ScrollRow = 2
$wsBusinessEvents.Range("A5").Select()

$wsScenarios.Select()
$wsScenarios.Range("A3").Select()
